$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 13 de Octubre de 2020 a las 00:33"

# Row 4
$ws.Range("B4").Value = 8034818
$ws.Range("C4").Value = 42820
$ws.Range("D4").Value = 5169007
$ws.Range("E4").Value = 2645813
$ws.Range("G4").Value = 303
$ws.Range("H4").Value = 219998

# Row 6
$ws.Range("D6").Value = 4495269
$ws.Range("E6").Value = 457450

# Row 8
$ws.Range("A8").Value = "Colombia"
$ws.Range("B8").Value = 919083
$ws.Range("C8").Value = 7767
$ws.Range("D8").Value = 798396
$ws.Range("E8").Value = 92702
$ws.Range("G8").Value = 151
$ws.Range("H8").Value = 27985

# Row 9
$ws.Range("A9").Value = "España"
$ws.Range("B9").Value = 918223
$ws.Range("C9").Value = 9286
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("G9").Value = 65
$ws.Range("H9").Value = 33124

# Row 46
$ws.Range("B46").Value = 104648
$ws.Range("C46").Value = 132
$ws.Range("D46").Value = 97743
$ws.Range("E46").Value = 843
$ws.Range("G46").Value = 10
$ws.Range("H46").Value = 6062

# Row 57
$ws.Range("B57").Value = 75948
$ws.Range("C57").Value = 334
$ws.Range("D57").Value = 71687
$ws.Range("E57").Value = 3981
$ws.Range("G57").Value = 5
$ws.Range("H57").Value = 280

# Row 69
$ws.Range("B69").Value = 47030
$ws.Range("C69").Value = 25
$ws.Range("D69").Value = 46424
$ws.Range("E69").Value = 298
$ws.Range("G69").Value = 2
$ws.Range("H69").Value = 308

# Row 85
$ws.Range("A85").Value = "Bulgaria"
$ws.Range("B85").Value = 24989
$ws.Range("C85").Value = 587
$ws.Range("D85").Value = 15975
$ws.Range("E85").Value = 8099
$ws.Range("G85").Value = 23
$ws.Range("H85").Value = 915

# Row 86
$ws.Range("A86").Value = "Corea del Sur"
$ws.Range("B86").Value = 24703
$ws.Range("C86").Value = 97
$ws.Range("D86").Value = 22729
$ws.Range("E86").Value = 1541
$ws.Range("G86").Value = 1
$ws.Range("H86").Value = 433

# Row 100
$ws.Range("B100").Value = 13691
$ws.Range("C100").Value = 6
$ws.Range("E100").Value = 6091

# Row 104
$ws.Range("B104").Value = 11062
$ws.Range("C104").Value = 40
$ws.Range("D104").Value = 10337
$ws.Range("E104").Value = 655
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 70

# Row 113
$ws.Range("B113").Value = 8860
$ws.Range("C113").Value = 25
$ws.Range("D113").Value = 8298
$ws.Range("E113").Value = 508

# Row 115
$ws.Range("B115").Value = 8021
$ws.Range("C115").Value = 10
$ws.Range("D115").Value = 7627
$ws.Range("E115").Value = 164

# Row 125
$ws.Range("B125").Value = 5426
$ws.Range("C125").Value = 3
$ws.Range("D125").Value = 5360
$ws.Range("E125").Value = 5

# Row 129
$ws.Range("B129").Value = 5116
$ws.Range("C129").Value = 15
$ws.Range("D129").Value = 3303
$ws.Range("E129").Value = 1721
$ws.Range("G129").Value = 2
$ws.Range("H129").Value = 92

# Row 131
$ws.Range("B131").Value = 5066
$ws.Range("C131").Value = 3
$ws.Range("D131").Value = 4954
$ws.Range("E131").Value = 29

# Row 160
$ws.Range("B160").Value = 2047
$ws.Range("C160").Value = 41
$ws.Range("E160").Value = 578

# Row 161
$ws.Range("B161").Value = 1949
$ws.Range("C161").Value = 9
$ws.Range("D161").Value = 1461
$ws.Range("E161").Value = 439

# Row 189
$ws.Range("B189").Value = 236
$ws.Range("C189").Value = 2
$ws.Range("D189").Value = 213
$ws.Range("E189").Value = 21
